# Revert "Add retry in case of server error"
#
# This undoes the addition of a "ServerErrorFailure" localization row
# (row 43) on the "Localization" sheet. Deleting the entire row shifts
# all subsequent rows up by one, automatically updates the sheet
# dimension, the backing table's ref/autoFilter, and drops the now
# unused shared-string entries ("ServerErrorFailure", its English
# description and its Japanese translation) from sharedStrings.xml.

$wb = $excel.ActiveWorkbook

$wsLocalization = $wb.Worksheets.Item("Localization")
$wsLocalization.Rows.Item(43).Delete()

# Restore "Settings" as the active/selected sheet (it was active before
# the original change moved the selection to "Localization").
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Activate()
